$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull of data
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = -3
$ws.Range("F17").Value = -5
